$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.636.76"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "2.550.70"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'301.79"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "'97.59"
$ws.Range("E6").Value = "  +6.68%  "
$ws.Range("D7").Value = "'0.573"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.544"
$ws.Range("E9").Value = "  -0.39%  "
$ws.Range("D10").Value = "'36.00"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").Value = "'0.0806"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("E12").Value = "  +8.72%  "
$ws.Range("D13").Value = "'7.49"
$ws.Range("E13").Value = "  -2.45%  "
$ws.Range("D14").Value = "2.539.41"
$ws.Range("E14").Value = "  -0.24%  "
$ws.Range("D15").Value = "'0.877"
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").Value = "'14.64"
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("D17").Value = "42.653.34"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "'13.29"
$ws.Range("E18").Value = "  +6.84%  "
$ws.Range("D19").Value = "0.0₃0981"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D20").Value = "'6.56"
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").Value = "'71.48"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Value = "'253.63"
$ws.Range("E22").Value = "  -2.44%  "
$ws.Range("D23").Value = "'2.93"
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("D25").Value = "'27.55"
$ws.Range("E25").Value = "  -6.67%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'10.01"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'37.75"
$ws.Range("E28").Value = "  +5.05%  "
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").Value = "'5.98"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").Value = "'154.85"
$ws.Range("E31").Value = "  +2.64%  "
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").Value = "'0.0801"
$ws.Range("E34").Value = "  +1.34%  "
$ws.Range("D35").Value = "'3.29"
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("B36").Value = "EnergySwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D36").Value = "'25.74"
$ws.Range("E36").Value = "  +5.85%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "'18.23"
$ws.Range("E37").Value = "  +13.29%  "
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("E40").Value = "  +32.14%  "
$ws.Range("D41").Value = "'3.86"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("D42").Value = "'3.36"
$ws.Range("E42").Value = "  -1.67%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0302"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.057.93"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("D46").Value = "'88.05"
$ws.Range("E46").Value = "  +3.61%  "
$ws.Range("D47").Value = "'9.20"
$ws.Range("E47").Value = "  +6.23%  "
$ws.Range("D48").Value = "2.796.40"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "'74.69"
$ws.Range("E49").Value = "  +8.10%  "
$ws.Range("D50").Value = "'102.95"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").Value = "'0.189"
$ws.Range("E51").Value = "  +1.39%  "
